# The sheet contains a daily price log for "Berenjena" (eggplant) at the
# "Vega Modelo de Temuco" market. A new daily record needs to be inserted
# at row 146 (it belongs there once the data is kept in date order along
# with the rest of the log), pushing every existing record from row 146
# down to row 147, and so on through to row 236 (what used to be the last
# record, on row 235, becomes row 236).
#
# Insert a whole new row at position 146 - Excel shifts rows 146:235 down
# to 147:236 automatically (and copies formatting, e.g. the date style on
# column D, from the row above), exactly mirroring what a user would do
# via right-click > "Insert" on the row header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("146:146").Insert()

# Populate the freshly inserted row with the new record's data.
$ws.Range("A146").Value = 10
$ws.Range("B146").Value = "Vega Modelo de Temuco"
$ws.Range("C146").Value = "La Araucanía"
$ws.Range("D146").Value = 44603
$ws.Range("E146").Value = 9
$ws.Range("F146").Value = 100112001
$ws.Range("G146").Value = "Berenjena"
$ws.Range("H146").Value = "Sin especificar"
$ws.Range("I146").Value = "Primera"
$ws.Range("J146").Value = 75
$ws.Range("K146").Value = 12000
$ws.Range("L146").Value = 13000
$ws.Range("M146").Value = 12467
$ws.Range("N146").Value = "$/caja 60 unidades"
$ws.Range("O146").Value = "Región del Maule"
$ws.Range("P146").Value = 208
$ws.Range("Q146").Value = 60
$ws.Range("R146").Value = "Hortaliza"
